# Implemented timeseries to supim file, demand file next
#
# 1) SupIm sheet: extend the yearly time series (t=1 template row) down
#    through t=12, i.e. rows 4..14, duplicating row 3's values/format.
# 2) Process sheet: merge the two conditional-formatting blocks
#    (A12:C13 A11 C11 / B11) into a single rule over A11:C13.
# 3) Make "SupIm" the active sheet/tab (was "Process").

$wb = $excel.ActiveWorkbook

# --- 1) SupIm: fill rows 4-14 with the same pattern as row 3 ---------------
$supim = $wb.Worksheets.Item("SupIm")

for ($r = 4; $r -le 14; $r++) {
    $t = $r - 2
    $supim.Cells.Item($r, 1).Value = $t
    $supim.Cells.Item($r, 2).Value = 0.481
    $supim.Cells.Item($r, 3).Value = 0.3
    $supim.Cells.Item($r, 4).Value = 0.207
}

# Copy formatting from the template row (row 3) down over the new rows.
$templateRow = $supim.Range("A3:D3")
$newRows = $supim.Range("A4:D14")
$templateRow.Copy()
$newRows.PasteSpecial(-4122)

# --- 2) Process sheet: consolidate conditional formatting ------------------
$process = $wb.Worksheets.Item("Process")
$fcs = $process.Cells.FormatConditions

# Remove the separate B11-only rule; keep the other rule and widen it to
# cover the whole A11:C13 block, making it first priority.
$fcB11 = $fcs.Item(2)
$fcB11.Delete()

$fcMain = $fcs.Item(1)
$fcMain.ModifyAppliesToRange($process.Range("A11:C13"))
$fcMain.SetFirstPriority()

# --- 3) Switch the active tab from Process to SupIm -------------------------
$supim.Activate() | Out-Null
$supim.Range("L16").Select() | Out-Null
